# Generate Report for Handoff
# Replaces the two "pending" sample files (a .png + a .md, each with a
# dependent .png) with two freshly-handed-off .md files, on all three
# sheets (Overview, zh-cn, de-de), and drops the now-obsolete
# "IsDependency" row from the zh-cn / de-de detail sheets.

$wb = $excel.ActiveWorkbook

$baseRepo = "https://github.com/OpenLocalizationTest/oltest/blob/6bf48540b0f53e99a2f8e388a6aa5334acf54d32"
$md1 = "126febc7-c600-460a-b705-01b79916ebbc.md"
$md2 = "52b0af62-2bc6-4228-bfcb-577d9923af7c.md"

# ---------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Drop all existing hyperlinks up front - they'll be re-created below
# once the final cell layout (after the row delete) is in place.
$ws1.Range("A1:C5").Hyperlinks.Delete()

# Row 4 (the second "*.png" dependency row) goes away entirely; row 5
# (".localization-config") slides up into row 4.
$ws1.Rows.Item(4).Delete()

$ws1.Range("A2").Value = $md1
$ws1.Range("A3").Value = $md2

$ws1.Hyperlinks.Add($ws1.Range("A2"), "$baseRepo/e2e/$md1", $null, $null, $md1)
$ws1.Hyperlinks.Add($ws1.Range("A3"), "$baseRepo/e2e/$md2", $null, $null, $md2)
$ws1.Hyperlinks.Add($ws1.Range("A4"), "$baseRepo/.localization-config", $null, $null, ".localization-config")

# ---------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A1:I5").Hyperlinks.Delete()

$xlf1zh = "126febc7-c600-460a-b705-01b79916ebbc.9e7cd770f8ec16970379c78675cdb15c13e87989.zh-cn.xlf"
$xlf2zh = "52b0af62-2bc6-4228-bfcb-577d9923af7c.ccd52850683e621c17fa1b2aacb1b885d32f80af.zh-cn.xlf"
$htBase1 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/93c9c4924ec336a49238972fcb1a3634419b75ae/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht"

# Row 4 (dependency row for the second png) is removed outright; row 5
# (".localization-config") becomes the new row 4.
$ws2.Rows.Item(4).Delete()

$ws2.Range("A2").Value = $md1
$ws2.Range("C2").Value = $xlf1zh
$ws2.Range("D2").Value = "2016-03-10 21:09:35"
$ws2.Range("H2").Value = "Include"
$ws2.Range("I2").ClearContents()

$ws2.Range("A3").Value = $md2
$ws2.Range("C3").Value = $xlf2zh
$ws2.Range("D3").Value = "2016-03-10 21:09:35"
$ws2.Range("H3").Value = "Include"

$ws2.Hyperlinks.Add($ws2.Range("A2"), "$baseRepo/e2e/$md1", $null, $null, $md1)
$ws2.Hyperlinks.Add($ws2.Range("C2"), "$htBase1/$xlf1zh", $null, $null, $xlf1zh)
$ws2.Hyperlinks.Add($ws2.Range("A3"), "$baseRepo/e2e/$md2", $null, $null, $md2)
$ws2.Hyperlinks.Add($ws2.Range("C3"), "$htBase1/$xlf2zh", $null, $null, $xlf2zh)
$ws2.Hyperlinks.Add($ws2.Range("A4"), "$baseRepo/.localization-config", $null, $null, ".localization-config")

# ---------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A1:I5").Hyperlinks.Delete()

$xlf1de = "126febc7-c600-460a-b705-01b79916ebbc.9e7cd770f8ec16970379c78675cdb15c13e87989.de-de.xlf"
$xlf2de = "52b0af62-2bc6-4228-bfcb-577d9923af7c.ccd52850683e621c17fa1b2aacb1b885d32f80af.de-de.xlf"
$htBase2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/42e0a2acccba3a2e20b920802a5cf898189f5ca0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht"

$ws3.Rows.Item(4).Delete()

$ws3.Range("A2").Value = $md1
$ws3.Range("C2").Value = $xlf1de
$ws3.Range("D2").Value = "2016-03-10 21:09:40"
$ws3.Range("H2").Value = "Include"
$ws3.Range("I2").ClearContents()

$ws3.Range("A3").Value = $md2
$ws3.Range("C3").Value = $xlf2de
$ws3.Range("D3").Value = "2016-03-10 21:09:40"
$ws3.Range("H3").Value = "Include"

$ws3.Hyperlinks.Add($ws3.Range("A2"), "$baseRepo/e2e/$md1", $null, $null, $md1)
$ws3.Hyperlinks.Add($ws3.Range("C2"), "$htBase2/$xlf1de", $null, $null, $xlf1de)
$ws3.Hyperlinks.Add($ws3.Range("A3"), "$baseRepo/e2e/$md2", $null, $null, $md2)
$ws3.Hyperlinks.Add($ws3.Range("C3"), "$htBase2/$xlf2de", $null, $null, $xlf2de)
$ws3.Hyperlinks.Add($ws3.Range("A4"), "$baseRepo/.localization-config", $null, $null, ".localization-config")
